# Correct the "Mean Age (Years)" row: replace the Unicode middle dot (·, U+00B7)
# used as a decimal separator with a standard period (.) so the values read
# as plain decimal numbers (e.g. "61·1 ± 12·7" -> "61.1 ± 12.7").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1")

$ws.Range("B2").Value = "61.1 ± 12.7"
$ws.Range("C2").Value = "59.8 ± 11.2"
$ws.Range("D2").Value = "60.6 ± 11.5"
$ws.Range("E2").Value = "60.5 ± 11.6"
